$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.002166666666666667
$ws.Range("H2").Value = 0.0065
$ws.Range("I2").Value = 0.004890446475191893
$ws.Range("J2").Value = 0.004890446475191893
$ws.Range("M2").Value = 9.084137666666667
$ws.Range("N2").Value = 27.252413
$ws.Range("O2").Value = 0.2765376761551382
$ws.Range("P2").Value = 0.2765376761551382
$ws.Range("Q2").Value = 0.01968229827777778
$ws.Range("R2").Value = 0.1771406845
$ws.Range("S2").Value = 0.001352392703610653
$ws.Range("T2").Value = 0.001352392703610653

# Row 3
$ws.Range("G3").Value = 0.002166666666666667
$ws.Range("H3").Value = 0.0065
$ws.Range("I3").Value = 0.004890446475191893
$ws.Range("J3").Value = 0.004890446475191893
$ws.Range("O3").Value = 0.3707916163717078
$ws.Range("P3").Value = 0.3707916163717078
$ws.Range("Q3").Value = 0.02639073016666667
$ws.Range("R3").Value = 0.2375165715
$ws.Range("S3").Value = 0.001813336553315723
$ws.Range("T3").Value = 0.001813336553315723

# Row 4
$ws.Range("G4").Value = 0.002166666666666667
$ws.Range("H4").Value = 0.0065
$ws.Range("I4").Value = 0.004890446475191893
$ws.Range("J4").Value = 0.004890446475191893
$ws.Range("M4").Value = 11.58507333333333
$ws.Range("N4").Value = 34.75522
$ws.Range("O4").Value = 0.3526707074731541
$ws.Range("P4").Value = 0.3526707074731541
$ws.Range("Q4").Value = 0.02510099222222222
$ws.Range("R4").Value = 0.22590893
$ws.Range("S4").Value = 0.001724717218265518
$ws.Range("T4").Value = 0.001724717218265518

# Row 5
$ws.Range("G5").Value = 0.440874
$ws.Range("H5").Value = 1.322622
$ws.Range("I5").Value = 0.9951095535248081
$ws.Range("J5").Value = 0.9951095535248081
$ws.Range("M5").Value = 9.084137666666667
$ws.Range("N5").Value = 27.252413
$ws.Range("O5").Value = 0.2765376761551382
$ws.Range("P5").Value = 0.2765376761551382
$ws.Range("Q5").Value = 4.004960109654
$ws.Range("R5").Value = 36.044640986886
$ws.Range("S5").Value = 0.2751852834515275
$ws.Range("T5").Value = 0.2751852834515275

# Row 6
$ws.Range("G6").Value = 0.440874
$ws.Range("H6").Value = 1.322622
$ws.Range("I6").Value = 0.9951095535248081
$ws.Range("J6").Value = 0.9951095535248081
$ws.Range("O6").Value = 0.3707916163717078
$ws.Range("P6").Value = 0.3707916163717078
$ws.Range("Q6").Value = 5.369993894538
$ws.Range("R6").Value = 48.32994505084201
$ws.Range("S6").Value = 0.3689782798183921
$ws.Range("T6").Value = 0.3689782798183921

# Row 7
$ws.Range("G7").Value = 0.440874
$ws.Range("H7").Value = 1.322622
$ws.Range("I7").Value = 0.9951095535248081
$ws.Range("J7").Value = 0.9951095535248081
$ws.Range("M7").Value = 11.58507333333333
$ws.Range("N7").Value = 34.75522
$ws.Range("O7").Value = 0.3526707074731541
$ws.Range("P7").Value = 0.3526707074731541
$ws.Range("Q7").Value = 5.10755762076
$ws.Range("R7").Value = 45.96801858684
$ws.Range("S7").Value = 0.3509459902548886
$ws.Range("T7").Value = 0.3509459902548886
